$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# D57 was stored as text "500043"; convert it to a real number.
$ws.Range("D57").Value = 500043

# Append six new rows (58-63) of stock data below the existing table.
# Columns: sr, nsecode, name, bsecode, per_chg, close, volume, timeframe, Date Time
$data = @(
    @(1, "ABBOTINDIA", "Abbott India Limited",        "500488", -0.95, 27377.45,       39754, "day", "27/06/2024 11:35:37"),
    @(2, "MARUTI",     "Maruti Suzuki India Limited",  "532500", -0.16, 12178.75,     1235718, "day", "27/06/2024 11:35:37"),
    @(3, "UBL",        "United Breweries Limited",     "532478",  1.35,  1990.70,      187019, "day", "27/06/2024 11:35:37"),
    @(4, "DALBHARAT",  "Dalmia Bharat Limited",        "533309", -1.52,  1792.20,     1110455, "day", "27/06/2024 11:35:37"),
    @(5, "TATAMOTORS", "Tata Motors Limited",          "500570",  2.13,   972.10,    19421905, "day", "27/06/2024 11:35:37"),
    @(6, "SUNTV",      "Sun Tv Network Limited",       "532733", -1.06,   754.75,     2923072, "day", "27/06/2024 11:35:37")
)

$row = 58
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]

    # bsecode (column D) must stay text, e.g. "500488", not become a number.
    # Force text storage via a temporary "@" format, then clear the format
    # back off so the written cell carries no style index (matches the
    # plain, unstyled inlineStr cells used throughout this sheet).
    $cellD = $ws.Cells.Item($row, 4)
    $cellD.NumberFormat = "@"
    $cellD.Value = $rec[3]
    $cellD.ClearFormats()

    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $ws.Cells.Item($row, 8).Value = $rec[7]
    $ws.Cells.Item($row, 9).Value = $rec[8]
    $row++
}
